$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 354.77777
$ws.Range("I32").Value = 257
$ws.Range("J32").Value = 403.66666
$ws.Range("K32").Value = 257
$ws.Range("L32").Value = 403.66666
$ws.Range("M32").Value = 69
$ws.Range("N32").Value = -1055.66666
$ws.Range("H47").Value = 17000
$ws.Range("I47").Value = 12000
$ws.Range("K47").Value = 12000
$ws.Range("M47").Value = -11028
$ws.Range("H64").Value = 3517.5952
$ws.Range("I64").Value = 3295.0435
$ws.Range("K64").Value = 3295.0435
$ws.Range("M64").Value = -3047.0435
$ws.Range("H67").Value = 3517.5952
$ws.Range("I67").Value = 3295.0435
$ws.Range("K67").Value = 3295.0435
$ws.Range("M67").Value = -2437.0435
$ws.Range("H98").Value = 1435.909
$ws.Range("I98").Value = 1649.375
$ws.Range("K98").Value = 1649.375
$ws.Range("M98").Value = -151.375
$ws.Range("H116").Value = 2664.4707
$ws.Range("I116").Value = 2436.3635
$ws.Range("J116").Value = 3082.6667
$ws.Range("K116").Value = 2436.3635
$ws.Range("L116").Value = 3082.6667
$ws.Range("M116").Value = 1005.6365
$ws.Range("N116").Value = -9966.6667
$ws.Range("H122").Value = 1435.909
$ws.Range("I122").Value = 1649.375
$ws.Range("K122").Value = 4948.125
$ws.Range("M122").Value = -2498.125
$ws.Range("H132").Value = 1970.1389
$ws.Range("I132").Value = 2197.111
$ws.Range("J132").Value = 1289.2222
$ws.Range("K132").Value = 6591.333
$ws.Range("L132").Value = 3867.6666
$ws.Range("M132").Value = -4061.333
$ws.Range("N132").Value = -8927.6666
$ws.Range("H137").Value = 3776.697
$ws.Range("I137").Value = 2308
$ws.Range("J137").Value = 5159
$ws.Range("K137").Value = 6924
$ws.Range("L137").Value = 15477
$ws.Range("M137").Value = -4374
$ws.Range("N137").Value = -20577
$ws.Range("H138").Value = 4735.1274
$ws.Range("I138").Value = 3709.5
$ws.Range("J138").Value = 5012.324
$ws.Range("K138").Value = 11128.5
$ws.Range("L138").Value = 15036.972
$ws.Range("M138").Value = -5988.5
$ws.Range("N138").Value = -25316.972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1535.6111
$ws.Range("I2").Value = 1874.3
$ws.Range("J2").Value = 1112.25
$ws.Range("K2").Value = 1874.3
$ws.Range("L2").Value = 1112.25
$ws.Range("M2").Value = -1761.3
$ws.Range("N2").Value = -1338.25
$ws.Range("H61").Value = 8536.471
$ws.Range("I61").Value = 6926.636
$ws.Range("J61").Value = 11487.833
$ws.Range("K61").Value = 6926.636
$ws.Range("L61").Value = 11487.833
$ws.Range("M61").Value = -6714.636
$ws.Range("N61").Value = -11911.833
$ws.Range("H97").Value = 967.6923
$ws.Range("I97").Value = 1031.6666
$ws.Range("J97").Value = 200
$ws.Range("K97").Value = 1031.6666
$ws.Range("L97").Value = 200
$ws.Range("M97").Value = -535.6666
$ws.Range("N97").Value = -1192
$ws.Range("H116").Value = 1535.6111
$ws.Range("I116").Value = 1874.3
$ws.Range("J116").Value = 1112.25
$ws.Range("K116").Value = 1874.3
$ws.Range("L116").Value = 1112.25
$ws.Range("M116").Value = 419.7
$ws.Range("N116").Value = -5700.25
$ws.Range("H136").Value = 8536.471
$ws.Range("I136").Value = 6926.636
$ws.Range("J136").Value = 11487.833
$ws.Range("K136").Value = 20779.908
$ws.Range("L136").Value = 34463.499
$ws.Range("M136").Value = -18229.908
$ws.Range("N136").Value = -39563.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1535.6111
$ws.Range("I3").Value = 1874.3
$ws.Range("J3").Value = 1112.25
$ws.Range("K3").Value = 1874.3
$ws.Range("L3").Value = 1112.25
$ws.Range("M3").Value = -1760.3
$ws.Range("N3").Value = -1340.25
$ws.Range("H86").Value = 3533.25
$ws.Range("I86").Value = 2740
$ws.Range("J86").Value = 7499.5
$ws.Range("K86").Value = 2740
$ws.Range("L86").Value = 7499.5
$ws.Range("M86").Value = -1617
$ws.Range("N86").Value = -9745.5
$ws.Range("H89").Value = 3533.25
$ws.Range("I89").Value = 2740
$ws.Range("J89").Value = 7499.5
$ws.Range("K89").Value = 13700
$ws.Range("L89").Value = 37497.5
$ws.Range("M89").Value = -8084
$ws.Range("N89").Value = -48729.5
$ws.Range("H94").Value = 1517.7222
$ws.Range("I94").Value = 944.1429000000001
$ws.Range("J94").Value = 1882.7273
$ws.Range("K94").Value = 944.1429000000001
$ws.Range("L94").Value = 1882.7273
$ws.Range("M94").Value = -493.1429000000001
$ws.Range("N94").Value = -2784.7273
$ws.Range("H118").Value = 57006.727
$ws.Range("J118").Value = 57006.727
$ws.Range("L118").Value = 57006.727
$ws.Range("N118").Value = -60320.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8963.315000000001
$ws.Range("I31").Value = 13599
$ws.Range("J31").Value = 4791.2
$ws.Range("K31").Value = 13599
$ws.Range("L31").Value = 4791.2
$ws.Range("M31").Value = -13304
$ws.Range("N31").Value = -5381.2
$ws.Range("H34").Value = 8963.315000000001
$ws.Range("I34").Value = 13599
$ws.Range("J34").Value = 4791.2
$ws.Range("K34").Value = 13599
$ws.Range("L34").Value = 4791.2
$ws.Range("M34").Value = -13397
$ws.Range("N34").Value = -5195.2
$ws.Range("H99").Value = 1885.069
$ws.Range("I99").Value = 1442.65
$ws.Range("J99").Value = 2868.2222
$ws.Range("K99").Value = 1442.65
$ws.Range("L99").Value = 2868.2222
$ws.Range("M99").Value = 55.34999999999991
$ws.Range("N99").Value = -5864.2222
$ws.Range("H109").Value = 49000
$ws.Range("J109").Value = 49000
$ws.Range("L109").Value = 49000
$ws.Range("N109").Value = -51080
$ws.Range("H126").Value = 1885.069
$ws.Range("I126").Value = 1442.65
$ws.Range("J126").Value = 2868.2222
$ws.Range("K126").Value = 4327.950000000001
$ws.Range("L126").Value = 8604.6666
$ws.Range("M126").Value = -1857.950000000001
$ws.Range("N126").Value = -13544.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4507847.5
$ws.Range("I5").Value = 416.26086
$ws.Range("J5").Value = 11912913
$ws.Range("K5").Value = 1248.78258
$ws.Range("L5").Value = 35738739
$ws.Range("M5").Value = -1136.78258
$ws.Range("N5").Value = -35738963
$ws.Range("H107").Value = 1272.8235
$ws.Range("I107").Value = 365.42856
$ws.Range("J107").Value = 1908
$ws.Range("K107").Value = 1096.28568
$ws.Range("L107").Value = 5724
$ws.Range("M107").Value = 823.71432
$ws.Range("N107").Value = -9564
$ws.Range("H110").Value = 4587.778
$ws.Range("I110").Value = 1450
$ws.Range("K110").Value = 4350
$ws.Range("M110").Value = -260
$ws.Range("H131").Value = 1314.8334
$ws.Range("I131").Value = 910.625
$ws.Range("J131").Value = 1516.9375
$ws.Range("K131").Value = 2731.875
$ws.Range("L131").Value = 4550.8125
$ws.Range("M131").Value = 2308.125
$ws.Range("N131").Value = -14630.8125
$ws.Range("H135").Value = 4507847.5
$ws.Range("I135").Value = 416.26086
$ws.Range("J135").Value = 11912913
$ws.Range("K135").Value = 3746.34774
$ws.Range("L135").Value = 107216217
$ws.Range("M135").Value = -1211.34774
$ws.Range("N135").Value = -107221287
$ws.Range("H136").Value = 3227.6428
$ws.Range("I136").Value = 1315
$ws.Range("J136").Value = 3992.7
$ws.Range("K136").Value = 3945
$ws.Range("L136").Value = 11978.1
$ws.Range("M136").Value = 1155
$ws.Range("N136").Value = -22178.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 44999.668
$ws.Range("J32").Value = 44999.668
$ws.Range("L32").Value = 44999.668
$ws.Range("N32").Value = -45591.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4430.9165
$ws.Range("I7").Value = 3383.923
$ws.Range("J7").Value = 5668.273
$ws.Range("K7").Value = 3383.923
$ws.Range("L7").Value = 5668.273
$ws.Range("M7").Value = -3271.923
$ws.Range("N7").Value = -5892.273
$ws.Range("H68").Value = 1580
$ws.Range("J68").Value = 2200
$ws.Range("L68").Value = 2200
$ws.Range("N68").Value = -3698
$ws.Range("H71").Value = 1580
$ws.Range("J71").Value = 2200
$ws.Range("L71").Value = 11000
$ws.Range("N71").Value = -18488
$ws.Range("H93").Value = 863.63635
$ws.Range("I93").Value = 900.8570999999999
$ws.Range("J93").Value = 798.5
$ws.Range("K93").Value = 900.8570999999999
$ws.Range("L93").Value = 798.5
$ws.Range("M93").Value = 347.1429000000001
$ws.Range("N93").Value = -3294.5
$ws.Range("H122").Value = 6764.222
$ws.Range("I122").Value = 6892.32
$ws.Range("J122").Value = 6473.091
$ws.Range("K122").Value = 20676.96
$ws.Range("L122").Value = 19419.273
$ws.Range("M122").Value = -18226.96
$ws.Range("N122").Value = -24319.273
$ws.Range("H126").Value = 4430.9165
$ws.Range("I126").Value = 3383.923
$ws.Range("J126").Value = 5668.273
$ws.Range("K126").Value = 10151.769
$ws.Range("L126").Value = 17004.819
$ws.Range("M126").Value = -7681.769
$ws.Range("N126").Value = -21944.819
$ws.Range("H132").Value = 3849.2563
$ws.Range("I132").Value = 3418.0688
$ws.Range("J132").Value = 5099.7
$ws.Range("K132").Value = 10254.2064
$ws.Range("L132").Value = 15299.1
$ws.Range("M132").Value = -7724.206399999999
$ws.Range("N132").Value = -20359.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 2000
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4746
$ws.Range("H109").Value = 58900
$ws.Range("J109").Value = 58900
$ws.Range("L109").Value = 58900
$ws.Range("N109").Value = -61674
$ws.Range("H122").Value = 3226.44
$ws.Range("I122").Value = 2582.8
$ws.Range("K122").Value = 7748.400000000001
$ws.Range("M122").Value = -5298.400000000001
$ws.Range("H132").Value = 2568.4443
$ws.Range("I132").Value = 1250
$ws.Range("J132").Value = 2945.1428
$ws.Range("K132").Value = 3750
$ws.Range("L132").Value = 8835.428400000001
$ws.Range("M132").Value = -1220
$ws.Range("N132").Value = -13895.4284
